# "Added correct data JFree Chart"
# Replace the old (incorrect) correlation results with the corrected
# Jacoco / Jfree Chart correlation data, re-format the sheet with a
# larger Times New Roman font, resize the columns and move the
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet - wipes stale values/formatting from the
# previous (wrong) data dump.
$ws.Cells.Clear()

# --- Column widths (character units read back from the COM layer do
# not map 1:1 onto the stored OOXML width, so these inputs are chosen
# to land on the desired stored widths of 58 / 57.1640625 / 28.5).
$ws.Columns("A").ColumnWidth = 57.166666666666664
$ws.Columns("B").ColumnWidth = 56.25
$ws.Columns("C").ColumnWidth = 27.666666666666668

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "Correlation"
$ws.Range("B1").Value = "Spearman Coefficient"
$ws.Range("C1").Value = "Type"

# --- Class level (Jacoco) correlations --------------------------------
$ws.Range("A3").Value = "Metrics 1-4  Jacoco version 1.0.19"
$ws.Range("B3").Value = 0.94516414444873198
$ws.Range("C3").Value = "Class Level"

$ws.Range("A4").Value = "Metrics  2-4  Jacoco version 1.0.19"
$ws.Range("B4").Value = 0.61341232088798803
$ws.Range("C4").Value = "Class Level"

$ws.Range("A6").Value = "Metrics 1-3  Jacoco version 1.0.19"
$ws.Range("B6").Value = 0.84179235354436199
$ws.Range("C6").Value = "Class Level"

$ws.Range("A7").Value = "Metrics 2-3  Jacoco version 1.0.19"
$ws.Range("B7").Value = 0.77145286098992305
$ws.Range("C7").Value = "Class Level"

# --- Version level (Jfree Chart) correlations -------------------------
$ws.Range("A9").Value = "Metrics 1-6  Jfree Chart 1.0.14,1.0.15,1.0.16,1.0.17,1.0.18,1.0.19"
$ws.Range("B9").Value = -0.98262379836822999
$ws.Range("C9").Value = "Version Level"

$ws.Range("A10").Value = "Metrics 2-6  Jfree Chart 1.0.14,1.0.15,1.0.16,1.0.17,1.0.18,1.0.19"
$ws.Range("B10").Value = -0.78262379212492605
$ws.Range("C10").Value = "Version Level"

$ws.Range("A12").Value = "Metrics  5-6  Jacoco version 1.0.19"
$ws.Range("B12").Value = -1
$ws.Range("C12").Value = "Version Level"

# --- Fonts --------------------------------------------------------------
# Base look for every populated cell: 16pt Times New Roman. Applied one
# (fully populated) row at a time so empty rows (2, 5, 8, 11) are never
# touched and stay absent from the sheet, matching the target layout.
foreach ($rowRange in @("A1:C1", "A3:C3", "A4:C4", "A6:C6", "A7:C7", "A9:C9", "A10:C10", "A12:C12")) {
    $r = $ws.Range($rowRange)
    $r.Font.Name = "Times New Roman"
    $r.Font.Size = 16
}

# Header row is bold.
$ws.Range("A1:B1").Font.Bold = $true

# The three "Class Level" coefficient cells carry an explicit black
# font color (rather than inheriting the automatic/theme color).
$ws.Range("B3").Font.Color = 0
$ws.Range("B4").Font.Color = 0
$ws.Range("B7").Font.Color = 0

# --- Selection ------------------------------------------------------
$ws.Range("B5").Select() | Out-Null
